# Apply crypto price/volume updates per the Feb 25 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.103.12'
$ws.Range('E2').Value = '  -3.55%  '
$ws.Range('D3').Value = '1.599.80'
$ws.Range('E3').Value = '  -2.93%  '
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '1.002'
$cell.Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  +0.04%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '300.89'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -2.89%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.3772'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -2.77%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.3627'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  -4.56%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '50.00'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  -3.72%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '1.250'
$cell.Style = 'Normal'
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '1.002'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +0.03%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.08115'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -3.77%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '22.92'
$cell.Style = 'Normal'
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '6.553'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -6.94%  '
$ws.Range('E15').Value = '  -3.99%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '7.348'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  -8.36%  '
$ws.Range('D17').Value = '1.599.88'
$ws.Range('E17').Value = '  -2.96%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '91.57'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  -2.63%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '0.06881'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -1.47%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '18.24'
$cell.Style = 'Normal'
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '6.532'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -5.75%  '
$ws.Range('B22').Value = 'BitDAO'
$ws.Range('C22').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '0.5567'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -6.82%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '1.002'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '12.93'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -5.76%  '
$ws.Range('B25').Value = 'WrappedBTC'
$ws.Range('C25').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D25').Value = '23.113.96'
$ws.Range('E25').Value = '  -3.46%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '2.338'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -4.77%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '2.714'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -7.89%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '21.02'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -4.58%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '149.61'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -2.50%  '
$ws.Range('B30').Value = 'HuobiToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '5.272'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -2.39%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '131.56'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -4.63%  '
$ws.Range('B32').Value = 'WEMIXTOKEN'
$ws.Range('C32').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '2.417'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -3.48%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '6.802'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -13.05%  '
$ws.Range('B34').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C34').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D34').Value = '1.774.50'
$ws.Range('E34').Value = '  -2.96%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '0.9520'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -6.24%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.07631'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -5.46%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '0.02717'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  -7.23%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.2536'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -5.06%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '6.185'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -8.16%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.08863'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -2.16%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '10.00'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -6.36%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '1.375'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -3.19%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '0.7048'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -6.87%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '12.58'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -5.80%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '15.35'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -5.26%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '0.6563'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  -5.39%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '2.287'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -6.29%  '
$ws.Range('B49').Value = 'PancakeSwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '3.973'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -2.91%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '131.55'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -1.76%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '0.07935'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -4.31%  '
